$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("programs")
$ws.Activate()

# Insert a new column before column B to hold the "load_data" source flag
# for each program parameter (yes/no whether it is loaded from external data).
$ws.Columns("B").Insert()

# Header row: rename the old "year" label to "program", and label the
# newly inserted column "load_data".
$ws.Range("A1").Value = "program"
$ws.Range("B1").Value = "load_data"

# Flag each program parameter row with whether it is currently wired up to
# loaded data. So far this is only true for vaccination.
$ws.Range("B2").Value = "yes"
$ws.Range("B3").Value = "no"
$ws.Range("B4").Value = "no"
$ws.Range("B5").Value = "no"
$ws.Range("B6").Value = "no"
$ws.Range("B7").Value = "no"

# Match the column width Excel would compute to fit "load_data".
$ws.Columns("B").ColumnWidth = 10.106445312500002

$ws.Range("G27").Select()
